$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for cells whose new values could otherwise be
# auto-detected by Excel as numbers (e.g. "0.999", "24.05"), then reset
# the cell style back to Normal so no spurious style/numFmt is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '64.629.75'
Set-TextValue $ws.Range('E2') '  -1.29%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.515.06'
Set-TextValue $ws.Range('E3') '  -1.77%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.01%  '

# Row 5
Set-TextValue $ws.Range('D5') '587.28'
Set-TextValue $ws.Range('E5') '  -2.54%  '

# Row 6
Set-TextValue $ws.Range('D6') '133.09'
Set-TextValue $ws.Range('E6') '  -2.01%  '

# Row 7
Set-TextValue $ws.Range('D7') '3.513.03'
Set-TextValue $ws.Range('E7') '  -1.75%  '

# Row 8
Set-TextValue $ws.Range('E8') '  +0.02%  '

# Row 9
Set-TextValue $ws.Range('E9') '  -1.07%  '

# Row 10
Set-TextValue $ws.Range('E10') '  +0.19%  '

# Row 11
Set-TextValue $ws.Range('D11') '7.15'
Set-TextValue $ws.Range('E11') '  -0.30%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.390'
Set-TextValue $ws.Range('E12') '  -0.35%  '

# Row 13
Set-TextValue $ws.Range('D13') '4.108.33'
Set-TextValue $ws.Range('E13') '  -1.83%  '

# Row 14
Set-TextValue $ws.Range('D14') '27.88'
Set-TextValue $ws.Range('E14') '  +1.21%  '

# Row 15
Set-TextValue $ws.Range('D15') '0.0000181'
Set-TextValue $ws.Range('E15') '  -2.35%  '

# Row 16
Set-TextValue $ws.Range('E16') '  +0.72%  '

# Row 17
Set-TextValue $ws.Range('D17') '3.510.37'
Set-TextValue $ws.Range('E17') '  -1.94%  '

# Row 18
Set-TextValue $ws.Range('D18') '64.598.33'
Set-TextValue $ws.Range('E18') '  -1.42%  '

# Row 19
Set-TextValue $ws.Range('D19') '9.98'
Set-TextValue $ws.Range('E19') '  -1.23%  '

# Row 20
Set-TextValue $ws.Range('D20') '14.29'
Set-TextValue $ws.Range('E20') '  -1.71%  '

# Row 21
Set-TextValue $ws.Range('E21') '  -2.78%  '

# Row 22
Set-TextValue $ws.Range('D22') '391.66'
Set-TextValue $ws.Range('E22') '  -0.40%  '

# Row 23
Set-TextValue $ws.Range('D23') '0.580'
Set-TextValue $ws.Range('E23') '  -0.56%  '

# Row 24
Set-TextValue $ws.Range('D24') '3.651.43'
Set-TextValue $ws.Range('E24') '  -1.85%  '

# Row 25
Set-TextValue $ws.Range('D25') '74.08'
Set-TextValue $ws.Range('E25') '  -0.17%  '

# Row 26
Set-TextValue $ws.Range('E26') '  +0.09%  '

# Row 27
Set-TextValue $ws.Range('E27') '  -3.17%  '

# Row 28
Set-TextValue $ws.Range('D28') '1.58'
Set-TextValue $ws.Range('E28') '  -7.51%  '

# Row 29
Set-TextValue $ws.Range('D29') '7.51'
Set-TextValue $ws.Range('E29') '  -7.05%  '

# Row 30
Set-TextValue $ws.Range('E30') '  +0.01%  '

# Row 31
Set-TextValue $ws.Range('E31') '  -2.91%  '

# Row 32
Set-TextValue $ws.Range('D32') '8.26'
Set-TextValue $ws.Range('E32') '  -4.62%  '

# Row 33
Set-TextValue $ws.Range('D33') '3.514.13'
Set-TextValue $ws.Range('E33') '  -1.90%  '

# Row 35
Set-TextValue $ws.Range('D35') '24.05'
Set-TextValue $ws.Range('E35') '  -1.18%  '

# Row 36
Set-TextValue $ws.Range('D36') '0.147'
Set-TextValue $ws.Range('E36') '  -0.99%  '

# Row 37
Set-TextValue $ws.Range('D37') '5.34'
Set-TextValue $ws.Range('E37') '  +2.63%  '

# Row 38
Set-TextValue $ws.Range('D38') '1.59'
Set-TextValue $ws.Range('E38') '  +1.07%  '

# Row 39
Set-TextValue $ws.Range('D39') '171.28'
Set-TextValue $ws.Range('E39') '  -0.27%  '

# Row 40
Set-TextValue $ws.Range('D40') '6.99'
Set-TextValue $ws.Range('E40') '  -0.40%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.0814'
Set-TextValue $ws.Range('E41') '  -1.98%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.814'
Set-TextValue $ws.Range('E42') '  -1.99%  '

# Row 43
Set-TextValue $ws.Range('D43') '26.27'
Set-TextValue $ws.Range('E43') '  -0.70%  '

# Row 44
Set-TextValue $ws.Range('E44') '  +0.05%  '

# Row 45
Set-TextValue $ws.Range('D45') '42.10'
Set-TextValue $ws.Range('E45') '  -2.44%  '

# Row 46
Set-TextValue $ws.Range('E46') '  -2.95%  '

# Row 47
Set-TextValue $ws.Range('D47') '4.42'
Set-TextValue $ws.Range('E47') '  -1.58%  '

# Row 48
Set-TextValue $ws.Range('E48') '  -2.18%  '

# Row 49
Set-TextValue $ws.Range('D49') '2.469.54'
Set-TextValue $ws.Range('E49') '  +0.15%  '

# Row 50
Set-TextValue $ws.Range('D50') '6.92'
Set-TextValue $ws.Range('E50') '  -0.91%  '

# Row 51
Set-TextValue $ws.Range('B51') 'SuiNetwork'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D51') '0.903'
Set-TextValue $ws.Range('E51') '  +2.88%  '
